# Applies the "Weighting and deadline" slide update:
#  - Overall weighting text: 25% -> 35%
#  - Drop the "12.5%" callouts from the two bullet sub-items
#  - Replace the big red deadline text (was "23rd May 2025 12 NOON")
#    with "25th March 2026, 12 NOON" as a single run, and resize/reposition
#    its text box
#  - Remove the now-unused "late deadline / no representation" notice textbox

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)

$nbsp = [char]0x00A0

# --- Shape 7: "TextBox 1" (weighting bullets) ---------------------------
$weightingShape = $s.Shapes.Item(7)
$weightingText = $weightingShape.TextFrame.TextRange

$weightingText.Paragraphs(1,1).Runs(1,1).Text = "Overall Assessment Weighting" + $nbsp + " (35% of 60 credits)"
$weightingText.Paragraphs(2,1).Runs(1,1).Text = "Part 1: Method Section (Turnitin)"
$weightingText.Paragraphs(3,1).Runs(1,1).Text = "Part 2: PsychoPy Experiment (Assignment)"

# --- Shape 8: "TextBox 4" (big red deadline date) ------------------------
$deadlineShape = $s.Shapes.Item(8)
$deadlineText = $deadlineShape.TextFrame.TextRange

# Originally 3 runs: "23" / "rd" (superscript) / " May 2025 12 NOON".
# Collapse them into a single run with the new wording.
$deadlineText.Runs(2,1).Text = ""
$deadlineText.Runs(1,1).Text = "25th March 2026, 12 NOON"
$deadlineText.Runs(2,1).Text = ""

$deadlineShape.Left = 191.52
$deadlineShape.Top = 347.5591339
$deadlineShape.Width = 602.64
$deadlineShape.Height = 60.58590551

# --- Shape 9: "TextBox 3" (late-deadline notice) -> remove --------------
$s.Shapes.Item(9).Delete()
